$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.481.43"
$ws.Range("E2").Value = "  -0.76%  "

# Row 3
$ws.Range("D3").Value = "1.825.42"
$ws.Range("E3").Value = "  -1.45%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "

# Row 6
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4259"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3613"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.70%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07208"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.21%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8637"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.57%  "

# Row 11
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.03%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.948.32"
$ws.Range("E12").Value = "  +6.52%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.385"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.01%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.475"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06944"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008935"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "

# Row 21
$ws.Range("D21").Value = "27.662.83"
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.124"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.57%  "

# Row 24
$ws.Range("D24").Value = "2.135.60"
$ws.Range("E24").Value = "  +2.74%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.992"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.135"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.793"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.39%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08901"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7476"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.985"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.543"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "

# Row 36
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
$ws.Range("E37").Value = "  -1.91%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05266"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.87%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01925"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.792"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5078"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

# Row 42
$ws.Range("E42").Value = "  -0.78%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.382"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.46%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.348"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.83%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.91%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.47"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4687"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.66%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06461"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.615"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "

Write-Output "Applied changes to cells"